$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "263.85"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "15"
$ws.Cells.Item(2, 7).Style = "Normal"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "22.65"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "15"
$ws.Cells.Item(3, 7).Style = "Normal"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "6.193"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "15"
$ws.Cells.Item(4, 7).Style = "Normal"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.06098"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "15"
$ws.Cells.Item(5, 7).Style = "Normal"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "3.532"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "15"
$ws.Cells.Item(6, 7).Style = "Normal"

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "6.727"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "15"
$ws.Cells.Item(7, 7).Style = "Normal"

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.364"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "15"
$ws.Cells.Item(8, 7).Style = "Normal"

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.8139"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "15"
$ws.Cells.Item(9, 7).Style = "Normal"

$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "15"
$ws.Cells.Item(10, 7).Style = "Normal"

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.08189"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "15"
$ws.Cells.Item(11, 7).Style = "Normal"

$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "15"
$ws.Cells.Item(12, 7).Style = "Normal"

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.03148"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "15"
$ws.Cells.Item(13, 7).Style = "Normal"

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.09260"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "15"
$ws.Cells.Item(14, 7).Style = "Normal"

$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "15"
$ws.Cells.Item(15, 7).Style = "Normal"

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.001707"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "15"
$ws.Cells.Item(16, 7).Style = "Normal"

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.04841"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "15"
$ws.Cells.Item(17, 7).Style = "Normal"

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.0006249"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "15"
$ws.Cells.Item(18, 7).Style = "Normal"

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.006258"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "15"
$ws.Cells.Item(19, 7).Style = "Normal"

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.006114"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "19HotbitTokenHTBBestin24h"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "15"
$ws.Cells.Item(20, 7).Style = "Normal"

$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "15"
$ws.Cells.Item(21, 7).Style = "Normal"

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.0001505"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "15"
$ws.Cells.Item(22, 7).Style = "Normal"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.695"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "15"
$ws.Cells.Item(23, 7).Style = "Normal"

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.266"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "15"
$ws.Cells.Item(24, 7).Style = "Normal"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.3404"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "15"
$ws.Cells.Item(25, 7).Style = "Normal"

$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "15"
$ws.Cells.Item(26, 7).Style = "Normal"

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0002689"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "15"
$ws.Cells.Item(27, 7).Style = "Normal"

$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "15"
$ws.Cells.Item(28, 7).Style = "Normal"

$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "15"
$ws.Cells.Item(29, 7).Style = "Normal"

$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "15"
$ws.Cells.Item(30, 7).Style = "Normal"

$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "15"
$ws.Cells.Item(31, 7).Style = "Normal"

$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "15"
$ws.Cells.Item(32, 7).Style = "Normal"

$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "15"
$ws.Cells.Item(33, 7).Style = "Normal"

$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "15"
$ws.Cells.Item(34, 7).Style = "Normal"

$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "15"
$ws.Cells.Item(35, 7).Style = "Normal"

$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "15"
$ws.Cells.Item(36, 7).Style = "Normal"

$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "15"
$ws.Cells.Item(37, 7).Style = "Normal"

$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "15"
$ws.Cells.Item(38, 7).Style = "Normal"

$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "15"
$ws.Cells.Item(39, 7).Style = "Normal"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.04643"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 7).NumberFormat = "@"
$ws.Cells.Item(40, 7).Value = "15"
$ws.Cells.Item(40, 7).Style = "Normal"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.007312"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 7).NumberFormat = "@"
$ws.Cells.Item(41, 7).Value = "15"
$ws.Cells.Item(41, 7).Style = "Normal"

$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1124"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"
$ws.Cells.Item(42, 7).NumberFormat = "@"
$ws.Cells.Item(42, 7).Value = "15"
$ws.Cells.Item(42, 7).Style = "Normal"

$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.003399"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "42CEJICEJI"
$ws.Cells.Item(43, 7).NumberFormat = "@"
$ws.Cells.Item(43, 7).Value = "15"
$ws.Cells.Item(43, 7).Style = "Normal"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.01036"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 7).NumberFormat = "@"
$ws.Cells.Item(44, 7).Value = "15"
$ws.Cells.Item(44, 7).Style = "Normal"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00006184"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 7).NumberFormat = "@"
$ws.Cells.Item(45, 7).Value = "15"
$ws.Cells.Item(45, 7).Style = "Normal"

$ws.Cells.Item(46, 7).NumberFormat = "@"
$ws.Cells.Item(46, 7).Value = "15"
$ws.Cells.Item(46, 7).Style = "Normal"

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.7522"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 7).NumberFormat = "@"
$ws.Cells.Item(47, 7).Value = "15"
$ws.Cells.Item(47, 7).Style = "Normal"

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.02614"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "47BOLOBOLO"
$ws.Cells.Item(48, 7).NumberFormat = "@"
$ws.Cells.Item(48, 7).Value = "15"
$ws.Cells.Item(48, 7).Style = "Normal"

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.00002107"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 7).NumberFormat = "@"
$ws.Cells.Item(49, 7).Value = "15"
$ws.Cells.Item(49, 7).Style = "Normal"

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.01244"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 7).NumberFormat = "@"
$ws.Cells.Item(50, 7).Value = "15"
$ws.Cells.Item(50, 7).Style = "Normal"

$ws.Cells.Item(51, 7).NumberFormat = "@"
$ws.Cells.Item(51, 7).Value = "15"
$ws.Cells.Item(51, 7).Style = "Normal"
